# Restructure ontology: remove mfd_hab1=Urban if mfd_areatype=Urban
#
# For every data row (2-83) in this P13_2 metadata sheet, mfd_areatype (L)
# is "Urban", so the redundant mfd_hab1="Urban" value is dropped: mfd_hab2
# ("Wastewater") takes its place in mfd_hab1, mfd_hab3 ("Activated sludge")
# takes the place of mfd_hab2, and mfd_hab3 is cleared (the ontology column
# list shrinks by one). habitat_typenumber (F) is also updated to the new
# code 2200.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 83

$rngF = $ws.Range("F$firstRow`:F$lastRow")
$rngN = $ws.Range("N$firstRow`:N$lastRow")
$rngO = $ws.Range("O$firstRow`:O$lastRow")
$rngP = $ws.Range("P$firstRow`:P$lastRow")

# habitat_typenumber: 1211 -> 2200, keep stored as text like the rest of
# the column (force text format so Excel doesn't coerce it to a number).
$rngF.NumberFormat = "@"
$rngF.Value = "2200"

# mfd_hab1 takes the old mfd_hab2 value ("Wastewater")
$rngN.Value = "Wastewater"

# mfd_hab2 takes the old mfd_hab3 value ("Activated sludge")
$rngO.Value = "Activated sludge"

# mfd_hab3 column is dropped for data rows
$rngP.ClearContents()
